$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Write-Host ("ActiveSheet: " + $ws.Name)
Write-Host ("B19: " + $ws.Range("B19").Value)
Write-Host ("H19: " + $ws.Range("H19").Value)
Write-Host ("I19: " + $ws.Range("I19").Value)
